# Refactor Code and Add Documentation
#
# The newest transaction row (row 2) originally held a placeholder /
# earlier "Bayar Cicilan" transfer. Replace it with the real
# "Top Up From BCA" transaction, and append the next top-up
# transaction as a brand-new row 3 (same data, later timestamp,
# running counter incremented).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    # Columns like Date ("2021-11-23") and Amount ("10000") look like a
    # date / a number to Excel's type-inference, so they'd otherwise be
    # silently coerced into a date-serial or numeric cell. Forcing the
    # cell to Text first keeps the literal string, then resetting the
    # style back to Normal drops the leftover "@" number-format so the
    # cell's formatting matches the rest of the (unstyled) data cells.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# --- Row 2: overwrite with the corrected transaction values ---
Set-TextValue $ws.Range("B2") "2021-11-23"
$ws.Range("C2").Value = "11:17:47"
$ws.Range("D2").Value = "Top Up From BCA"
$ws.Range("E2").Value = "Top Up"
Set-TextValue $ws.Range("F2") "10000"
$ws.Range("G2").Value = "BCA"
$ws.Range("H2").Value = "MichaelH"

# --- Row 3: new transaction row, appended after row 2 ---
$ws.Range("A3").Value = 1
Set-TextValue $ws.Range("B3") "2021-11-23"
$ws.Range("C3").Value = "11:18:03"
$ws.Range("D3").Value = "Top Up From BCA"
$ws.Range("E3").Value = "Top Up"
Set-TextValue $ws.Range("F3") "10000"
$ws.Range("G3").Value = "BCA"
$ws.Range("H3").Value = "MichaelH"

# A3 (the running counter cell) should carry the same bordered /
# centered style as A2, so copy A2's formatting across.
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$excel.CutCopyMode = $false
